$wb = $excel.ActiveWorkbook
$newSheet = $wb.Worksheets.Add()
$newSheet.Name = "PDP"

$newSheet.Range("A1:Y1").Interior.Color = 65535

$newSheet.Range("A1").Value = "DataSet"
$newSheet.Range("B1").Value = "UserName"
$newSheet.Range("C1").Value = "Password"
$newSheet.Range("D1").Value = "Confirm Password"
$newSheet.Range("E1").Value = "FirstName"
$newSheet.Range("F1").Value = "LastName"
$newSheet.Range("G1").Value = "Email"
$newSheet.Range("H1").Value = "methods"
$newSheet.Range("I1").Value = "Street"
$newSheet.Range("J1").Value = "City"
$newSheet.Range("K1").Value = "Region"
$newSheet.Range("L1").Value = "postcode"
$newSheet.Range("M1").Value = "phone"
$newSheet.Range("N1").Value = "OTP Number"
$newSheet.Range("O1").Value = "cardNumber"
$newSheet.Range("P1").Value = "ExpMonthYear"
$newSheet.Range("Q1").Value = "cvv"
$newSheet.Range("R1").Value = "Products"
$newSheet.Range("S1").Value = "Colorproduct"
$newSheet.Range("T1").Value = "Color"
$newSheet.Range("U1").Value = "Quantity"
$newSheet.Range("V1").Value = "Discountcode"
$newSheet.Range("W1").Value = "DOB"
$newSheet.Range("X1").Value = "Links"
$newSheet.Range("Y1").Value = "productquantity"

$newSheet.Range("A2").Value = "Product"
$newSheet.Range("R2").Value = "32 oz Wide Mouth "
$newSheet.Range("S2").Value = "32 oz Wide Mouth"
$newSheet.Range("T2").Value = "Black"
$newSheet.Range("Y2").Value = 2

$newSheet.Range("Y2").Select()
